# Update "想去人数" (F column) and "最低票价" (G column) figures on both the
# "展览" and "全部类型" sheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 521
    $ws.Range("G4").Value = 50
    $ws.Range("F8").Value = 495
    $ws.Range("F9").Value = 6539
    $ws.Range("F11").Value = 137
    $ws.Range("F12").Value = 1024
    $ws.Range("F13").Value = 344
    $ws.Range("F14").Value = 108
    $ws.Range("F15").Value = 180
    $ws.Range("F16").Value = 492
}
